$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New 2x2 layout: header row (first_name/second_name) + one data row (f1/l1),
# replacing the old 1x4 layout (title/firstname/lastname/company).
$ws.Range("A1").Value = "first_name"
$ws.Range("B1").Value = "second_name"
$ws.Range("A2").Value = "f1"
$ws.Range("B2").Value = "l1"

# The old C1/D1 cells are no longer part of the used range.
$ws.Range("C1").Value = $null
$ws.Range("D1").Value = $null

# Touch the formatting of the new range so it carries its own cell style
# (applies the "Normal" style explicitly instead of leaving cells on the
# sheet's implicit default style).
$ws.Range("A1:B2").Style = "Normal"

# Move the active selection the way the authored workbook has it.
$ws.Range("C2").Select() | Out-Null
